# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Ají" (Cristal, Macroferia Regional de
# Talca) as row 346, pushing the existing rows 346-361 down to 347-362.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 346 - this shifts rows 346..361 down to
# 347..362 and extends the used range / dimension to A1:R362 automatically.
$ws.Rows("346").Insert()

# Populate the new row 346 with the new weekly record.
$ws.Cells.Item(346, 1).Value = 5
$ws.Cells.Item(346, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(346, 3).Value = "Maule"
$ws.Cells.Item(346, 4).Value = 45041
$ws.Cells.Item(346, 5).Value = 7
$ws.Cells.Item(346, 6).Value = 100112021
$ws.Cells.Item(346, 7).Value = "Ají"
$ws.Cells.Item(346, 8).Value = "Cristal"
$ws.Cells.Item(346, 9).Value = "Primera"
$ws.Cells.Item(346, 10).Value = 100
$ws.Cells.Item(346, 11).Value = 15000
$ws.Cells.Item(346, 12).Value = 15000
$ws.Cells.Item(346, 13).Value = 15000
$ws.Cells.Item(346, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(346, 15).Value = "Región del Maule"
$ws.Cells.Item(346, 16).Value = 600
$ws.Cells.Item(346, 17).Value = 25
$ws.Cells.Item(346, 18).Value = "Hortaliza"
